$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "coldread_gaze_wpm_std" row (row 6) first so the row index of the
# "total_duration" row (row 3) is unaffected by this deletion.
$ws.Rows.Item(6).Delete()

# Remove the "total_duration" row; everything below shifts up one.
$ws.Rows.Item(3).Delete()

# Fix the typo in the coverage description (missing closing parenthesis after
# "title") - this cell is now at C3 after the row deletions above.
$ws.Range("C3").Value = "The coverage of X (word, line, paragraph, title) by fixation in percentage during the coldread section"

# Append the new qa_saccade_regression_rate_% feature as the new last row.
$ws.Range("A8").Value = "qa_saccade_regression_rate_%"
$ws.Range("B8").Value = "Strategy"
$ws.Range("C8").Value = "The percentage of saccades that were regression within the text."

# Update the view: scroll back to the top-left (A1) and select B4.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
